# 3 Scenes translated - Lime's fj, Lily's temptation, and liliy's thiighjob
#
# The source sheet has Japanese text in column A and (where translated)
# the English text in column D. This adds the four missing English
# translations for the "locked room" flavor text entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "メイドの部屋のようだ。鍵がかかっている・・・" (Housemaid's room)
$ws.Range("D10").Value = "This seems to be the Housemaid's room.`nIt's locked tightly..."

# Row 11: "扉が開かない・・・魔法的な力で閉じられているようだ。" (The door won't open)
$ws.Range("D11").Value = "The door won't open...`nIt seems to be closed shut by a magical power."

# Row 29: "書庫のようだ。鍵がかかっている・・・" (Library)
$ws.Range("D29").Value = "This seems to be the Library.`nIt's locked tightly..."

# Row 31: "館主の部屋のようだ。鍵がかかっている・・・" (household owner's office)
$ws.Range("D31").Value = "This seems to be office of the household owner.`nIt's locked tightly..."

# Re-fit the affected rows so embedding a line-break in the new text doesn't
# leave behind an explicit row-height override (matches how the rest of the
# sheet's multi-line cells are left at the default row height).
$ws.Range("D10").EntireRow.AutoFit()
$ws.Range("D11").EntireRow.AutoFit()
$ws.Range("D29").EntireRow.AutoFit()
$ws.Range("D31").EntireRow.AutoFit()
